# Update CDA Logical model for ST.r2b
# - Bump Version / Date metadata values
# - Insert a new "Jurisdiction" property row (blank value) before "Description"

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Insert a new row at position 11 (pushes Description.. down by one row)
$ws.Rows.Item(11).Insert()

# Copy the formatting of the row above (a regular data row) onto the new row
# so the new cells get the same style as the rest of the property rows.
$ws.Range("A10:B10").Copy()
$ws.Range("A11:B11").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Populate the new "Jurisdiction" row
$ws.Range("A11").Value = "Jurisdiction"
$ws.Range("B11").Value = ""

# Update the Version value
$ws.Range("B3").Value = "2.0.1-sd-202510-matchbox-patch"

# Update the Date value
$ws.Range("B8").Value = "2025-10-29T22:15:57+01:00"
